# Apply the crypto-list price/volume refresh captured by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.993.44"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "3.422.11"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'572.53"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'163.73"
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.422.05"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "'0.555"
$ws.Range("E9").Value = "  -4.75%  "
$ws.Range("D10").Value = "'7.30"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").Value = "'0.120"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").Value = "'0.425"
$ws.Range("E12").Value = "  -4.76%  "
$ws.Range("D13").Value = "4.016.23"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").Value = "'27.09"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").Value = "'0.0000174"
$ws.Range("E16").Value = "  -6.33%  "
$ws.Range("D17").Value = "64.050.85"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "3.375.10"
$ws.Range("E18").Value = "  -0.87%  "
$ws.Range("D19").Value = "'6.15"
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").Value = "'13.69"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").Value = "'378.06"
$ws.Range("E21").Value = "  -1.09%  "
$ws.Range("D22").Value = "'7.80"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'71.20"
$ws.Range("E24").Value = "  -1.40%  "
$ws.Range("D25").Value = "'0.518"
$ws.Range("E25").Value = "  -5.72%  "
$ws.Range("E26").Value = "  -1.31%  "
$ws.Range("D27").Value = "'9.56"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "'6.11"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = "  -5.55%  "
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").Value = "'22.98"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "'7.14"
$ws.Range("E34").Value = "  +1.62%  "
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("D36").Value = "'160.02"
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").Value = "'0.858"
$ws.Range("E37").Value = "  +10.28%  "
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0729"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.800.92"
$ws.Range("E40").Value = "  -3.77%  "
$ws.Range("D41").Value = "'25.90"
$ws.Range("E41").Value = "  -1.43%  "
$ws.Range("D42").Value = "'42.84"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("D44").Value = "'26.36"
$ws.Range("E44").Value = "  +1.39%  "
$ws.Range("D45").Value = "'4.42"
$ws.Range("E45").Value = "  -3.09%  "
$ws.Range("D46").Value = "'0.0307"
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("D47").Value = "'2.45"
$ws.Range("E47").Value = "  +8.47%  "
$ws.Range("D48").Value = "'330.09"
$ws.Range("E48").Value = "  +4.04%  "
$ws.Range("E49").Value = "  -4.23%  "
$ws.Range("D50").Value = "'6.34"
$ws.Range("E50").Value = "  -2.89%  "
$ws.Range("E51").Value = "  -2.49%  "

# A handful of the Price values above are plain decimals (e.g. "7.30", "0.998") and were
# written with a leading apostrophe so Excel keeps them as text (matching the other Price
# cells, which are already text because of the thousands-dot formatting, e.g. "63.993.44").
# That quote-prefix marks the cells with a "number stored as text" style, so re-apply the
# plain Normal style across the whole data range to keep formatting identical to the rest
# of the sheet (the data cells carry no explicit style in the original workbook).
$ws.Range("B2:E51").Style = "Normal"
